$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6504219999999999
$ws.Range("H2").Value = 1.951266
$ws.Range("I2").Value = 0.1521898546336546
$ws.Range("J2").Value = 0.1521898546336546
$ws.Range("M2").Value = 1.01111
$ws.Range("N2").Value = 3.03333
$ws.Range("O2").Value = 0.04063212692754557
$ws.Range("P2").Value = 0.04063212692754556
$ws.Range("Q2").Value = 0.6576481884199999
$ws.Range("R2").Value = 5.918833695779999
$ws.Range("S2").Value = 0.006183797490559363
$ws.Range("T2").Value = 0.006183797490559362

# Row 3
$ws.Range("G3").Value = 0.6504219999999999
$ws.Range("H3").Value = 1.951266
$ws.Range("I3").Value = 0.1521898546336546
$ws.Range("J3").Value = 0.1521898546336546
$ws.Range("O3").Value = 0.4065982422683317
$ws.Range("P3").Value = 0.4065982422683317
$ws.Range("Q3").Value = 6.580964809431332
$ws.Range("R3").Value = 59.228683284882
$ws.Range("S3").Value = 0.06188012738511689
$ws.Range("T3").Value = 0.06188012738511689

# Row 4
$ws.Range("G4").Value = 0.6504219999999999
$ws.Range("H4").Value = 1.951266
$ws.Range("I4").Value = 0.1521898546336546
$ws.Range("J4").Value = 0.1521898546336546
$ws.Range("O4").Value = 0.5527696308041227
$ws.Range("P4").Value = 0.5527696308041226
$ws.Range("Q4").Value = 8.946810659460667
$ws.Range("R4").Value = 80.521295935146
$ws.Range("S4").Value = 0.08412592975797836
$ws.Range("T4").Value = 0.08412592975797835

# Row 5
$ws.Range("I5").Value = 0.437647089654669
$ws.Range("J5").Value = 0.4376470896546689
$ws.Range("M5").Value = 1.01111
$ws.Range("N5").Value = 3.03333
$ws.Range("O5").Value = 0.04063212692754557
$ws.Range("P5").Value = 0.04063212692754556
$ws.Range("Q5").Value = 1.89117609956
$ws.Range("R5").Value = 17.02058489604
$ws.Range("S5").Value = 0.01778253209631942
$ws.Range("T5").Value = 0.01778253209631942

# Row 6
$ws.Range("I6").Value = 0.437647089654669
$ws.Range("J6").Value = 0.4376470896546689
$ws.Range("O6").Value = 0.4065982422683317
$ws.Range("P6").Value = 0.4065982422683317
$ws.Range("S6").Value = 0.1779465373874394
$ws.Range("T6").Value = 0.1779465373874394

# Row 7
$ws.Range("I7").Value = 0.437647089654669
$ws.Range("J7").Value = 0.4376470896546689
$ws.Range("O7").Value = 0.5527696308041227
$ws.Range("P7").Value = 0.5527696308041226
$ws.Range("S7").Value = 0.2419180201709101
$ws.Range("T7").Value = 0.2419180201709101

# Row 8
$ws.Range("I8").Value = 0.4101630557116764
$ws.Range("J8").Value = 0.4101630557116764
$ws.Range("M8").Value = 1.01111
$ws.Range("N8").Value = 3.03333
$ws.Range("O8").Value = 0.04063212692754557
$ws.Range("P8").Value = 0.04063212692754556
$ws.Range("Q8").Value = 1.77241111896
$ws.Range("R8").Value = 15.95170007064
$ws.Range("S8").Value = 0.01666579734066678
$ws.Range("T8").Value = 0.01666579734066678

# Row 9
$ws.Range("I9").Value = 0.4101630557116764
$ws.Range("J9").Value = 0.4101630557116764
$ws.Range("O9").Value = 0.4065982422683317
$ws.Range("P9").Value = 0.4065982422683317
$ws.Range("S9").Value = 0.1667715774957755
$ws.Range("T9").Value = 0.1667715774957755

# Row 10
$ws.Range("I10").Value = 0.4101630557116764
$ws.Range("J10").Value = 0.4101630557116764
$ws.Range("O10").Value = 0.5527696308041227
$ws.Range("P10").Value = 0.5527696308041226
$ws.Range("S10").Value = 0.2267256808752342
$ws.Range("T10").Value = 0.2267256808752341
